# "added common method for search"
# Adds four new location rows (loc1-loc4) to the "add_new_locations" sheet,
# flips the "runmode" flag (Y/N) for the existing SLIIT / IDM rows, and
# nudges the workbook window width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- flip runmode (column D) for existing rows 2 and 4 (Y -> N) ---
$ws.Range("D2").Value = "N"
$ws.Range("D4").Value = "N"

# --- new rows, written column-by-column so shared strings land in the
#     same order Excel produced them in (B, then C, then A, then D) ---
$ws.Range("B5").Value = "Battaramulla"
$ws.Range("B6").Value = "Kandy"
$ws.Range("B7").Value = "Kalutara"
$ws.Range("B8").Value = "Jaela"

$ws.Range("C5").Value = "Active"
$ws.Range("C6").Value = "Active"
$ws.Range("C7").Value = "Inactive"
$ws.Range("C8").Value = "Inactive"

$ws.Range("A5").Value = "loc1"
$ws.Range("A6").Value = "loc2"
$ws.Range("A7").Value = "loc3"
$ws.Range("A8").Value = "loc4"

$ws.Range("D5").Value = "Y"
$ws.Range("D6").Value = "Y"
$ws.Range("D7").Value = "Y"
$ws.Range("D8").Value = "Y"

# --- window / selection bookkeeping to mirror the committed workbook ---
$excel.ActiveWindow.Width = 4935

$null = $ws.Range("A9").Select()
